# Data format fixes for CSV
#
# The "Tasks" sheet gains a new "finish" column (a 0/1 flag, matching the
# pattern already used on the "Subs"/"Subsubs" sheets) positioned right
# before the existing "expired" column. Concretely: insert a new column at
# E, shifting the old E column ("expired" / "3 days") to F, then populate
# the new E column with the header "finish" and the value 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Insert a new blank column at E; everything from E onward (old "expired"
# header/value) shifts right to F.
$ws.Columns("E:E").Insert()

# Populate the freshly inserted column E with the new "finish" data.
$ws.Range("E1").Value = "finish"
$ws.Range("E2").Value = 1
